# working on Gen start log
$wb = $excel.ActiveWorkbook

$todo  = $wb.Worksheets.Item("TODO")
$login = $wb.Worksheets.Item("登录")

# remember TODO's new selection (does not change the active sheet)
$todo.Range("D20").Select()

# --- 1. Add the new worksheet "机组启停记录" at the end ---
$lastSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "机组启停记录"
$newSheet.Move($null, $wb.Worksheets.Item($lastSheetName))
$ws = $wb.Worksheets.Item("机组启停记录")

# --- 2. Fill in the "机组启停记录" sheet content ---
# (ordered to mirror how the data was typed in originally)
$ws.Range("C9").Value = "页面"
$ws.Range("D4").Value = "机组编号"
$ws.Range("E4").Value = "启动时间"
$ws.Range("F4").Value = "记录人"
$ws.Range("G4").Value = "停止时间"
$ws.Range("C3").Value = "db"
$ws.Range("C4").Value = "电厂编号"
$ws.Range("C16").Value = "操作"
$ws.Range("C17").Value = "添加"
$ws.Range("D12").Value = "修改"
$ws.Range("E12").Value = "删除"
$ws.Range("D18").Value = "修改、删除，限制录入人操作"
$ws.Range("C10").Value = "page-1"
$ws.Range("C11").Value = "page-2"
$ws.Range("D11").Value = "查询"
$ws.Range("F12").Value = "导出"
$ws.Range("F13").Value = "excel文件"
$ws.Range("F5").Value = "中文姓名"
$ws.Range("C14").Value = "page-3"
$ws.Range("D14").Value = "图标"
$ws.Range("I4").Value = "停止标志"
$ws.Range("I5").Value = "1-未停止"
$ws.Range("I6").Value = "2-已停止"

# cells that reuse already-registered strings
$ws.Range("H4").Value = "记录人"
$ws.Range("D10").Value = "机组编号"
$ws.Range("E10").Value = "启动时间"
$ws.Range("F10").Value = "停止时间"
$ws.Range("D17").Value = "修改"
$ws.Range("E17").Value = "删除"

# date in A3, formatted like the date cells on the TODO sheet
$todo.Range("B2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 43577

# red font styling (copied from the existing red-font cell on 登录!G7)
$login.Range("G7").Copy()
$ws.Range("F5:F8").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)

$ws.Columns.Item(1).ColumnWidth = 9.2857142857
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- 3. Add the new row of data to the "TODO" sheet ---
$todo.Range("A3").Value = 2
$todo.Range("B2").Copy()
$todo.Range("B3").PasteSpecial(-4122)
$todo.Range("B3").Value = 43577
$todo.Range("C3").Value = "Login后 store用户哪些信息"

# --- 4. Leave "机组启停记录" as the active sheet/tab ---
$ws.Activate()
$ws.Range("H21").Select()
